$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues - paste back as a literal value (not a live formula)
$xlPasteValues = -4163

$dCell = $ws.Range("D2")
$dCell.Formula = '="96.828.37"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E2").Value = "  +0.47%  "
$dCell = $ws.Range("D3")
$dCell.Formula = '="3.663.76"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$dCell = $ws.Range("D5")
$dCell.Formula = '="242.80"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  +20.55%  "
$dCell = $ws.Range("D7")
$dCell.Formula = '="657.28"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("E10").Value = "  -0.08%  "
$dCell = $ws.Range("D11")
$dCell.Formula = '="3.661.73"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E11").Value = "  +2.27%  "
$dCell = $ws.Range("D12")
$dCell.Formula = '="44.34"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E12").Value = "  +2.64%  "
$dCell = $ws.Range("D13")
$dCell.Formula = '="0.204"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +1.92%  "
$dCell = $ws.Range("D15")
$dCell.Formula = '="4.342.87"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E15").Value = "  +2.19%  "
$dCell = $ws.Range("D16")
$dCell.Formula = '="96.584.19"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E16").Value = "  +0.29%  "
$dCell = $ws.Range("D17")
$dCell.Formula = '="0.0000259"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E17").Value = "  -0.32%  "
$dCell = $ws.Range("D18")
$dCell.Formula = '="3.636.54"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E18").Value = "  +1.60%  "
$dCell = $ws.Range("D19")
$dCell.Formula = '="8.18"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("E20").Value = "  +3.26%  "
$dCell = $ws.Range("D21")
$dCell.Formula = '="18.37"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E21").Value = "  +3.61%  "
$dCell = $ws.Range("D22")
$dCell.Formula = '="0.529"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E22").Value = "  +6.89%  "
$dCell = $ws.Range("D23")
$dCell.Formula = '="512.59"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E23").Value = "  +0.18%  "
$dCell = $ws.Range("D24")
$dCell.Formula = '="3.44"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  +0.58%  "
$dCell = $ws.Range("D27")
$dCell.Formula = '="101.44"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E27").Value = "  +5.21%  "
$dCell = $ws.Range("D28")
$dCell.Formula = '="13.07"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("E29").Value = "  +13.77%  "
$ws.Range("E30").Value = "  +1.46%  "
$dCell = $ws.Range("D31")
$dCell.Formula = '="11.86"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E31").Value = "  +3.65%  "
$dCell = $ws.Range("D32")
$dCell.Formula = '="1.00"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +1.95%  "
$dCell = $ws.Range("D34")
$dCell.Formula = '="33.21"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E34").Value = "  +5.11%  "
$dCell = $ws.Range("D35")
$dCell.Formula = '="0.997"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("E37").Value = "  +3.68%  "
$dCell = $ws.Range("D38")
$dCell.Formula = '="615.89"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +1.41%  "
$dCell = $ws.Range("D40")
$dCell.Formula = '="42.45"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E40").Value = "  +24.30%  "
$ws.Range("E41").Value = "  +5.80%  "
$dCell = $ws.Range("D42")
$dCell.Formula = '="0.955"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E42").Value = "  +5.34%  "
$dCell = $ws.Range("D43")
$dCell.Formula = '="1.94"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E43").Value = "  +6.15%  "
$ws.Range("E44").Value = "  +0.01%  "
$dCell = $ws.Range("D45")
$dCell.Formula = '="6.14"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E45").Value = "  +7.75%  "
$dCell = $ws.Range("D46")
$dCell.Formula = '="0.0443"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E46").Value = "  +6.17%  "
$dCell = $ws.Range("D47")
$dCell.Formula = '="0.423"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E47").Value = "  +26.93%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +0.43%  "
$dCell = $ws.Range("D50")
$dCell.Formula = '="8.61"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E50").Value = "  +5.35%  "
$dCell = $ws.Range("D51")
$dCell.Formula = '="54.56"'
$dCell.Copy()
$dCell.PasteSpecial($xlPasteValues)
$ws.Range("E51").Value = "  +2.97%  "
